# penghapusan tombol edit dan hapus, penambahan pembagian zonasi di datawajibpajak
#
# The sheet "Sheet1" held a small sample/demo table:
#   A1: npwpd              B1: jumlah_penagihan
#   A2: P.2.0021240.03.003 B2: 120000
#   A3: P.2.0021556.01.011 B3: 2000
#   A4: P.2.0001580.04.009.B4: 400000
#
# This template is consumed by the app at render time (rows get filled in
# programmatically per-wajib-pajak / zona), so the hard-coded sample rows
# are removed here, leaving just the header row and three blank,
# still-styled rows ready to be populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three sample data rows (npwpd + jumlah_penagihan values) while
# keeping row/cell formatting (style) intact so the template still looks
# right once new data is written back in.
$ws.Range("A2:B4").ClearContents()

# The "npwpd" column used a slightly different font (Helvetica Neue) than
# the rest of the template (Arial) -- unify on Arial.
$ws.Range("A2:A4").Font.Name = "Arial"

# Leave the selection parked over the now-empty data rows, matching where
# an editor would click next to start filling the template back in.
$ws.Range("A2:B4").Select()
